$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text storage for numeric-looking Price values so they stay strings (matches original inlineStr cells)
$textCells = @("D5", "D6", "D7", "D8", "D9", "D10", "D13", "D14", "D15", "D16", "D20", "D22", "D23", "D24", "D25", "D26", "D27", "D28", "D29", "D30", "D31", "D32", "D34", "D37", "D38", "D39", "D40", "D41", "D42", "D43", "D44", "D45", "D46", "D47", "D48", "D49", "D50", "D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply updated values
$ws.Range("D2").Value = "30.485.82"
$ws.Range("E2").Value = "  -0.35%  "
$ws.Range("D3").Value = "1.908.97"
$ws.Range("E3").Value = "  -0.62%  "
$ws.Range("D5").Value = "243.67"
$ws.Range("E5").Value = "  -1.71%  "
$ws.Range("D6").Value = "1.000"
$ws.Range("E6").Value = "  +0.02%  "
$ws.Range("D7").Value = "0.4835"
$ws.Range("D8").Value = "0.2882"
$ws.Range("E8").Value = "  -0.32%  "
$ws.Range("D9").Value = "0.06799"
$ws.Range("E9").Value = "  -0.50%  "
$ws.Range("D10").Value = "111.45"
$ws.Range("E10").Value = "  +6.06%  "
$ws.Range("E11").Value = "  +5.10%  "
$ws.Range("D12").Value = "1.918.71"
$ws.Range("E12").Value = "  -0.04%  "
$ws.Range("D13").Value = "0.07559"
$ws.Range("E13").Value = "  -1.84%  "
$ws.Range("D14").Value = "5.380"
$ws.Range("E14").Value = "  +1.58%  "
$ws.Range("D15").Value = "0.6669"
$ws.Range("E15").Value = "  -0.29%  "
$ws.Range("D16").Value = "292.42"
$ws.Range("D17").Value = "30.498.25"
$ws.Range("E17").Value = "  -0.32%  "
$ws.Range("E18").Value = "  +0.30%  "
$ws.Range("E19").Value = "  +0.09%  "
$ws.Range("D20").Value = "0.000007555"
$ws.Range("E20").Value = "  -0.61%  "
$ws.Range("D21").Value = "2.163.73"
$ws.Range("E21").Value = "  -0.21%  "
$ws.Range("D22").Value = "5.482"
$ws.Range("E22").Value = "  -1.09%  "
$ws.Range("D23").Value = "1.000"
$ws.Range("D24").Value = "6.377"
$ws.Range("E24").Value = "  -0.03%  "
$ws.Range("D25").Value = "9.426"
$ws.Range("E25").Value = "  +0.20%  "
$ws.Range("D26").Value = "165.26"
$ws.Range("E26").Value = "  -1.60%  "
$ws.Range("D27").Value = "20.18"
$ws.Range("E27").Value = "  -4.38%  "
$ws.Range("D28").Value = "2.064"
$ws.Range("E28").Value = "  -2.27%  "
$ws.Range("D29").Value = "0.1064"
$ws.Range("E29").Value = "  -0.37%  "
$ws.Range("D30").Value = "1.428"
$ws.Range("E30").Value = "  +2.21%  "
$ws.Range("D31").Value = "4.122"
$ws.Range("E31").Value = "  -1.34%  "
$ws.Range("D32").Value = "4.046"
$ws.Range("E32").Value = "  -0.58%  "
$ws.Range("E33").Value = "  -1.20%  "
$ws.Range("D34").Value = "0.7338"
$ws.Range("E34").Value = "  -0.67%  "
$ws.Range("E35").Value = "  -1.32%  "
$ws.Range("B37").Value = "HuobiToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D37").Value = "2.718"
$ws.Range("E37").Value = "  -0.92%  "
$ws.Range("B38").Value = "VeChain"
$ws.Range("C38").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D38").Value = "0.02030"
$ws.Range("E38").Value = "  -2.14%  "
$ws.Range("D39").Value = "2.680"
$ws.Range("E39").Value = "  -0.31%  "
$ws.Range("B40").Value = "Quant"
$ws.Range("C40").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D40").Value = "109.23"
$ws.Range("E40").Value = "  -1.66%  "
$ws.Range("B41").Value = "RenderToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D41").Value = "2.004"
$ws.Range("E41").Value = "  -2.40%  "
$ws.Range("D42").Value = "0.4412"
$ws.Range("E42").Value = "  +0.73%  "
$ws.Range("D43").Value = "0.8606"
$ws.Range("E43").Value = "  -1.93%  "
$ws.Range("D44").Value = "5.769"
$ws.Range("E44").Value = "  -1.70%  "
$ws.Range("D45").Value = "1.0000"
$ws.Range("E45").Value = "  +0.03%  "
$ws.Range("D46").Value = "68.97"
$ws.Range("E46").Value = "  +1.88%  "
$ws.Range("D47").Value = "7.163"
$ws.Range("E47").Value = "  -1.26%  "
$ws.Range("D48").Value = "48.10"
$ws.Range("E48").Value = "  -0.68%  "
$ws.Range("D49").Value = "9.215"
$ws.Range("E49").Value = "  -0.87%  "
$ws.Range("D50").Value = "0.1224"
$ws.Range("E50").Value = "  -0.66%  "
$ws.Range("D51").Value = "0.2511"
$ws.Range("E51").Value = "  +0.66%  "
